$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.823.36"
$ws.Range("E2").Value = "  +4.30%  "
$ws.Range("D3").Value = "2.275.46"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.41%  "
$ws.Range("E6").Value = "  +5.81%  "
$ws.Range("E7").Value = "  +3.80%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +3.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0800"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("D15").Value = "2.627.53"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").Value = "2.283.54"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("E18").Value = "  +3.65%  "
$ws.Range("D19").Value = "41.767.52"
$ws.Range("E19").Value = "  +4.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.03%  "
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("E25").Value = "  +4.99%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("E35").Value = "  +4.62%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.05%  "
$ws.Range("E39").Value = "  +5.48%  "
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("E41").Value = "  +3.76%  "
$ws.Range("E42").Value = "  +5.53%  "
$ws.Range("D43").Value = "2.077.15"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("E45").Value = "  +3.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("E47").Value = "  +7.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.41%  "
$ws.Range("E51").Value = "  +3.55%  "
